$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes
# Note: Excel COM's ColumnWidth setter stores width using a Max-Digit-Width based
# rounding, so to land exactly on an integer stored width we nudge the requested
# value slightly (empirically verified against this workbook's font/styles).
$ws.Columns("M").ColumnWidth = 19.15
$ws.Columns("P").ColumnWidth = 21.15
$ws.Columns("V").ColumnWidth = 19.15

$ws.Range("M15").Value = -2.65
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 3.75
$ws.Range("P15").Value = -0.15
$ws.Range("W15").Value = 0

$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0.959
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 10.441
$ws.Range("W16").Value = 0

$ws.Range("M17").Value = 49.996
$ws.Range("N17").Value = 0.959
$ws.Range("O17").Value = -46.396
$ws.Range("P17").Value = -0.3089999999999999

$ws.Range("M18").Value = 47
$ws.Range("N18").Value = -37.541
$ws.Range("O18").Value = -47
$ws.Range("P18").Value = 47.84099999999999

$ws.Range("M19").Value = 49.996
$ws.Range("N19").Value = -26.537
$ws.Range("O19").Value = -58.746
$ws.Range("P19").Value = 26.537

$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 6.4
$ws.Range("P20").Value = 0
$ws.Range("W20").Value = 0

$ws.Range("M21").Value = -44
$ws.Range("N21").Value = 0.959
$ws.Range("O21").Value = 34.65
$ws.Range("P21").Value = -0.959

$ws.Range("M22").Value = 44
$ws.Range("N22").Value = 0.959
$ws.Range("O22").Value = -53.35
$ws.Range("P22").Value = -0.959

$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = -12.35

$ws.Range("M24").Value = -49.95
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 59.2
$ws.Range("P24").Value = 0

$ws.Range("M25").Value = -0.3
$ws.Range("N25").Value = -0.767
$ws.Range("O25").Value = 2.05
$ws.Range("P25").Value = 1.167
$ws.Range("W25").Value = 0

$ws.Range("M26").Value = -0.5629999999999999
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 7.613
$ws.Range("P26").Value = -2.15
$ws.Range("W26").Value = 0

$ws.Range("N27").Value = -50.562
$ws.Range("P27").Value = 63.11199999999999

$ws.Range("M28").Value = -50.562
$ws.Range("O28").Value = 39.962

$ws.Range("M93").Value = 0.1369999999999999
$ws.Range("N93").Value = 0
$ws.Range("O93").Value = -1.387
$ws.Range("P93").Value = 0
$ws.Range("W93").Value = 0

$ws.Range("M94").Value = 1.52
$ws.Range("N94").Value = 0
$ws.Range("O94").Value = -0.52
$ws.Range("P94").Value = 0.4
$ws.Range("W94").Value = 0

$ws.Range("M95").Value = 0.958
$ws.Range("N95").Value = 0
$ws.Range("O95").Value = -0.958
$ws.Range("P95").Value = 0
$ws.Range("W95").Value = 0

$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = -11.25

$ws.Range("M97").Value = -49.95
$ws.Range("N97").Value = 0
$ws.Range("O97").Value = 60.3
$ws.Range("P97").Value = 0

$ws.Range("M98").Value = 49.996
$ws.Range("N98").Value = 0.959
$ws.Range("O98").Value = -46.396
$ws.Range("P98").Value = -1.309

$ws.Range("M99").Value = -0.3
$ws.Range("N99").Value = 0.756
$ws.Range("O99").Value = 3.1
$ws.Range("P99").Value = -0.106
$ws.Range("W99").Value = 0

$ws.Range("M100").Value = -0.3
$ws.Range("N100").Value = -0.662
$ws.Range("O100").Value = 3.1
$ws.Range("P100").Value = 0.412
$ws.Range("W100").Value = 0

$ws.Range("M101").Value = 0
$ws.Range("N101").Value = 0.959
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = 11.041
$ws.Range("W101").Value = 0

$ws.Range("M102").Value = 47
$ws.Range("N102").Value = -37.541
$ws.Range("O102").Value = -47
$ws.Range("P102").Value = 48.34099999999999

$ws.Range("M103").Value = 49.996
$ws.Range("N103").Value = -26.537
$ws.Range("O103").Value = -58.146
$ws.Range("P103").Value = 27.337

$ws.Range("M104").Value = -44
$ws.Range("N104").Value = 0.959
$ws.Range("O104").Value = 34.05
$ws.Range("P104").Value = -0.959

$ws.Range("M105").Value = 44
$ws.Range("N105").Value = 0.959
$ws.Range("O105").Value = -53.95
$ws.Range("P105").Value = -0.959

